$d = $word.ActiveDocument

# The edit strips the document down to just the "This is a Markdown file"
# Heading2 paragraph: the Title/Author/Date paragraphs above it, and every
# paragraph below it (R Markdown boilerplate prose, code chunks, the
# "Including Plots" section, the plot image, etc.) are removed.

$keepText = "This is a Markdown file"
$keepParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq $keepText) {
        $keepParagraph = $p
        break
    }
}

if ($keepParagraph -eq $null) {
    throw "Could not locate the '$keepText' paragraph"
}

$keepStart = $keepParagraph.Range.Start
$keepEnd = $keepParagraph.Range.End

# Delete the tail first so the offsets for the head deletion stay valid.
$docEnd = $d.Content.End
if ($docEnd -gt $keepEnd) {
    $tail = $d.Range($keepEnd, $docEnd)
    $tail.Delete()
}

if ($keepStart -gt 0) {
    $head = $d.Range(0, $keepStart)
    $head.Delete()
}
